$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 61:62. Row 60 (with the original "Primera" Ciruela
# record for 2022-03-02 / O'Higgins) stays put for now; the row that used to
# be 61 (the "Segunda" record) is pushed down to row 63.
$ws.Rows("61:62").Insert()

# Row 62 now becomes a duplicate of the (still unmodified) original row 60
# content, since row 60 itself is about to be overwritten with new values.
$ws.Range("A62").Value2 = 7
$ws.Range("B62").Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C62").Value2 = 'Ñuble'
$ws.Range("D62").Value2 = 44622
$ws.Range("E62").Value2 = 16
$ws.Range("F62").Value2 = 'Fruta'
$ws.Range("G62").Value2 = 100103
$ws.Range("H62").Value2 = 'Frutos de hueso (carozo)'
$ws.Range("I62").Value2 = 100103002
$ws.Range("J62").Value2 = 'Ciruela'
$ws.Range("K62").Value2 = 'Black Amber'
$ws.Range("L62").Value2 = 'Primera'
$ws.Range("M62").Value2 = 240
$ws.Range("N62").Value2 = 11000
$ws.Range("O62").Value2 = 12000
$ws.Range("P62").Value2 = 11500
$ws.Range("Q62").Value2 = '$/bandeja 18 kilos granel'
$ws.Range("R62").Value2 = "Región de O'Higgins"
$ws.Range("S62").Value2 = 639
$ws.Range("T62").Value2 = 18

# Row 61 is a brand-new record: "Segunda" Ciruela Black Amber, 2023-01-05,
# Provincia de Curicó.
$ws.Range("A61").Value2 = 7
$ws.Range("B61").Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C61").Value2 = 'Ñuble'
$ws.Range("D61").Value2 = 44931
$ws.Range("E61").Value2 = 16
$ws.Range("F61").Value2 = 'Fruta'
$ws.Range("G61").Value2 = 100103
$ws.Range("H61").Value2 = 'Frutos de hueso (carozo)'
$ws.Range("I61").Value2 = 100103002
$ws.Range("J61").Value2 = 'Ciruela'
$ws.Range("K61").Value2 = 'Black Amber'
$ws.Range("L61").Value2 = 'Segunda'
$ws.Range("M61").Value2 = 60
$ws.Range("N61").Value2 = 13000
$ws.Range("O61").Value2 = 13000
$ws.Range("P61").Value2 = 13000
$ws.Range("Q61").Value2 = '$/bandeja 18 kilos granel'
$ws.Range("R61").Value2 = 'Provincia de Curicó'
$ws.Range("S61").Value2 = 722
$ws.Range("T61").Value2 = 18

# Finally, update row 60 in place: "Primera" Ciruela Black Amber record moves
# to 2023-01-05 / Provincia de Curicó with new prices.
$ws.Range("D60").Value2 = 44931
$ws.Range("M60").Value2 = 120
$ws.Range("N60").Value2 = 15000
$ws.Range("O60").Value2 = 16000
$ws.Range("P60").Value2 = 15500
$ws.Range("R60").Value2 = 'Provincia de Curicó'
$ws.Range("S60").Value2 = 861
